$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 242
$ws.Range("F3").Value = 1420
$ws.Range("F5").Value = 893
$ws.Range("F6").Value = 31
$ws.Range("F7").Value = 1250
$ws.Range("F8").Value = 1598
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 2280
$ws.Range("F12").Value = 459
$ws.Range("F13").Value = 125
$ws.Range("F14").Value = 47
$ws.Range("F16").Value = 95
$ws.Range("F17").Value = 84
$ws.Range("F18").Value = 6257
$ws.Range("F20").Value = 6156
$ws.Range("F21").Value = 10150
$ws.Range("F23").Value = 176
$ws.Range("F25").Value = 279
$ws.Range("F26").Value = 503
$ws.Range("F27").Value = 177
$ws.Range("F28").Value = 154
$ws.Range("F29").Value = 4397
$ws.Range("F30").Value = 88
$ws.Range("F31").Value = 391

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 1160

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 242
$ws.Range("F5").Value = 1420
$ws.Range("F8").Value = 893
$ws.Range("F9").Value = 31
$ws.Range("F10").Value = 1250
$ws.Range("F12").Value = 1598
$ws.Range("F15").Value = 2280
$ws.Range("F17").Value = 459
$ws.Range("F18").Value = 125
$ws.Range("F19").Value = 47
$ws.Range("F22").Value = 95
$ws.Range("F23").Value = 84
$ws.Range("F24").Value = 6258
$ws.Range("F26").Value = 6156
$ws.Range("F27").Value = 10150
$ws.Range("F30").Value = 176
$ws.Range("F32").Value = 279
$ws.Range("F34").Value = 503
$ws.Range("F38").Value = 177
$ws.Range("F39").Value = 154
$ws.Range("F40").Value = 4397
$ws.Range("F46").Value = 391
